# Updated symbol list on Tue Feb 14 05:49:16 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the crypto
# tracker sheet with the latest pulled quotes. Values are written as text
# (matching the sheet's existing inline-string cells for these columns),
# so the target cells are pre-formatted as Text before the write to stop
# Excel from auto-converting the numeric-looking price strings / percent
# strings into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "292.21" },
    @{ Cell = "E2";  Value = "-7.25%" },
    @{ Cell = "D3";  Value = "40.31" },
    @{ Cell = "E3";  Value = "-1.55%" },
    @{ Cell = "D4";  Value = "5.031" },
    @{ Cell = "E4";  Value = "-2.50%" },
    @{ Cell = "D5";  Value = "0.07323" },
    @{ Cell = "E5";  Value = "-3.56%" },
    @{ Cell = "D6";  Value = "1.534" },
    @{ Cell = "E6";  Value = "-7.88%" },
    @{ Cell = "D7";  Value = "0.9274" },
    @{ Cell = "E7";  Value = "0.01%" },
    @{ Cell = "D9";  Value = "0.1186" },
    @{ Cell = "E9";  Value = "-1.11%" },
    @{ Cell = "D10"; Value = "0.1745" },
    @{ Cell = "E10"; Value = "-4.08%" },
    @{ Cell = "E11"; Value = "4.05%" },
    @{ Cell = "D12"; Value = "0.08654" },
    @{ Cell = "E12"; Value = "-4.31%" },
    @{ Cell = "D13"; Value = "0.1055" },
    @{ Cell = "E13"; Value = "0.24%" },
    @{ Cell = "D14"; Value = "0.001268" },
    @{ Cell = "E14"; Value = "-1.11%" },
    @{ Cell = "D15"; Value = "0.006005" },
    @{ Cell = "E15"; Value = "3.51%" },
    @{ Cell = "D16"; Value = "3.339" },
    @{ Cell = "E16"; Value = "0.25%" },
    @{ Cell = "D17"; Value = "4.295" },
    @{ Cell = "D19"; Value = "7.977" },
    @{ Cell = "E19"; Value = "5.41%" },
    @{ Cell = "D20"; Value = "0.1400" },
    @{ Cell = "E20"; Value = "3.59%" },
    @{ Cell = "D21"; Value = "0.2744" },
    @{ Cell = "E21"; Value = "-2.12%" },
    @{ Cell = "E22"; Value = "-2.50%" },
    @{ Cell = "D23"; Value = "0.001262" },
    @{ Cell = "E23"; Value = "-0.85%" },
    @{ Cell = "D24"; Value = "0.003780" },
    @{ Cell = "E24"; Value = "-6.70%" },
    @{ Cell = "D25"; Value = "0.0001282" },
    @{ Cell = "E25"; Value = "0.93%" },
    @{ Cell = "D26"; Value = "0.0003727" },
    @{ Cell = "D38"; Value = "0.02276" },
    @{ Cell = "E38"; Value = "-6.01%" },
    @{ Cell = "D39"; Value = "0.04978" },
    @{ Cell = "E39"; Value = "-3.59%" },
    @{ Cell = "E40"; Value = "70.69%" },
    @{ Cell = "D41"; Value = "0.007698" },
    @{ Cell = "E41"; Value = "-0.31%" },
    @{ Cell = "D42"; Value = "0.1285" },
    @{ Cell = "E42"; Value = "-1.23%" },
    @{ Cell = "D43"; Value = "0.007329" },
    @{ Cell = "E43"; Value = "-3.64%" },
    @{ Cell = "D44"; Value = "0.008291" },
    @{ Cell = "E44"; Value = "-3.09%" },
    @{ Cell = "D45"; Value = "0.2915" },
    @{ Cell = "E45"; Value = "-14.45%" },
    @{ Cell = "D46"; Value = "0.00006306" },
    @{ Cell = "E46"; Value = "-4.15%" },
    @{ Cell = "E47"; Value = "0.14%" },
    @{ Cell = "D48"; Value = "0.02559" },
    @{ Cell = "E48"; Value = "-90.70%" },
    @{ Cell = "D49"; Value = "0.00002103" },
    @{ Cell = "E49"; Value = "0.14%" },
    @{ Cell = "D50"; Value = "0.0002003" },
    @{ Cell = "E50"; Value = "0.14%" }
)

foreach ($update in $updates) {
    $rng = $ws.Range($update.Cell)
    # Force text storage (these columns hold numeric-looking / percent-looking
    # strings as plain text, not numbers) so the write doesn't get silently
    # reinterpreted as a Number/Percentage by Excel's smart-entry parsing.
    $rng.NumberFormat = "@"
    $rng.Value = $update.Value
}

Write-Host "Updated $($updates.Count) cells"
